$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 3 (Ängsmetallvinge / Adscita statices record) is being pushed
# down to become row 4, and row 3 is replaced with a new record
# (Svartfläckig blåvinge / Phengaris arion) that has several extra columns filled in.

# Insert a new blank row at position 4, then duplicate the (still unmodified)
# row 3 into it so the old record is preserved as row 4.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0

# Now overwrite row 3 with the new record's data.
$ws.Range("A3").Value = 112144581
$ws.Range("B3").Value = 42594
$ws.Range("E3").Value = 101260
$ws.Range("F3").Value = "Svartfläckig blåvinge"
$ws.Range("G3").Value = "Phengaris arion"

# "Antal" (I3) is stored as text ("1"), not a number, in the source data,
# so force text formatting before assigning it, then drop the formatting
# override again (it is only needed to make Excel keep the value as text).
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"
$ws.Range("I3").ClearFormats()

$ws.Range("J3").Value = "ex."
$ws.Range("K3").Value = "imago/adult"
$ws.Range("L3").Value = "hona"
$ws.Range("M3").Value = "vilande"

$ws.Range("Q3").Value = 442971.9404393921
$ws.Range("R3").Value = 6204766.971186478

$ws.Range("AC3").Value = "lufthåvning"

$ws.Range("AO3").Value = "på grässtrå"
$ws.Range("AQ3").Value = "Nils Otto Nilsson"
$ws.Range("AR3").Value = "NON 04616"
